# "added css from branch feature-c"
#
# Append a new paragraph at the very end of the document containing the
# text "add feature from feature c - branch", matching the look of the
# existing paragraph (Calibri 11pt body text, default spacing/indent,
# left justified).

$d = $word.ActiveDocument

# Collapse a range positioned at the end of the document's main story so
# the new paragraph is added after everything that is already there.
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)   # wdCollapseEnd

$newPara = $d.Paragraphs.Add($endOfDoc)
$newPara.Range.Text = "add feature from feature c - branch"

# Match the font used throughout the rest of the document.
$newPara.Range.Font.Name = "Calibri"
$newPara.Range.Font.Size = 11
